# Apply cryptos list update (prices & 1h volume % changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a number by Excel
# (single-decimal values like "580.47" or "2.00") must be forced to Text format
# first so the literal string survives exactly as scraped.
$textForceRows = @(5,6,9,10,11,12,13,14,19,20,21,22,23,24,25,29,30,31,32,33,35,36,37,39,40,43,46,48,51)
foreach ($r in $textForceRows) {
    $ws.Range("D" + $r).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.687.14"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "3.035.27"
$ws.Range("E3").Value = "  -4.47%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "580.47"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "129.67"
$ws.Range("E6").Value = "  -4.23%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.035.59"
$ws.Range("E8").Value = "  -4.41%  "
$ws.Range("D9").Value = "0.501"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("D11").Value = "5.22"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").Value = "33.41"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "3.537.36"
$ws.Range("E16").Value = "  -4.41%  "
$ws.Range("D17").Value = "61.757.81"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "3.037.17"
$ws.Range("E18").Value = "  -4.37%  "
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("D20").Value = "443.94"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").Value = "13.37"
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("D22").Value = "0.667"
$ws.Range("E22").Value = "  -5.00%  "
$ws.Range("D23").Value = "7.28"
$ws.Range("E23").Value = "  -4.41%  "
$ws.Range("D24").Value = "80.36"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("D25").Value = "12.77"
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").Value = "2.00"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").Value = "7.43"
$ws.Range("E30").Value = "  -4.16%  "
$ws.Range("D31").Value = "6.43"
$ws.Range("E31").Value = "  -5.84%  "
$ws.Range("D32").Value = "25.82"
$ws.Range("E32").Value = "  -5.31%  "
$ws.Range("D33").Value = "0.0965"
$ws.Range("E33").Value = "  -6.80%  "
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("D35").Value = "0.972"
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("D36").Value = "5.67"
$ws.Range("E36").Value = "  -4.08%  "
$ws.Range("D37").Value = "50.28"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").Value = "0.0₃0699"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "0.0371"
$ws.Range("E39").Value = "  -3.53%  "
$ws.Range("D40").Value = "7.87"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("E42").Value = "  -7.18%  "
$ws.Range("D43").Value = "376.84"
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("D44").Value = "2.674.69"
$ws.Range("E44").Value = "  -4.46%  "
$ws.Range("D46").Value = "122.67"
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("E47").Value = "  -4.75%  "
$ws.Range("D48").Value = "34.09"
$ws.Range("E48").Value = "  -5.93%  "
$ws.Range("E49").Value = "  -6.17%  "
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "23.69"
$ws.Range("E51").Value = "  -6.38%  "
